# Re-run SGNN to annotate dialog acts following clean up work to the
# original transcripts. This updates the DAMSLTag (column I) and
# DialogAct (column J) values for a set of rows in the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 18;  Tag = "sd";  Act = "Statement-non-opinion" }
    @{ Row = 21;  Tag = "%";   Act = "Uninterpretable" }
    @{ Row = 32;  Tag = "aa";  Act = "Agree/Accept" }
    @{ Row = 34;  Tag = "sd";  Act = "Statement-non-opinion" }
    @{ Row = 42;  Tag = "aa";  Act = "Agree/Accept" }
    @{ Row = 45;  Tag = "aa";  Act = "Agree/Accept" }
    @{ Row = 48;  Tag = "aa";  Act = "Agree/Accept" }
    @{ Row = 55;  Tag = "aa";  Act = "Agree/Accept" }
    @{ Row = 57;  Tag = "aa";  Act = "Agree/Accept" }
    @{ Row = 61;  Tag = "aa";  Act = "Agree/Accept" }
    @{ Row = 69;  Tag = "sd";  Act = "Statement-non-opinion" }
    @{ Row = 73;  Tag = "sd";  Act = "Statement-non-opinion" }
    @{ Row = 76;  Tag = "%";   Act = "Uninterpretable" }
    @{ Row = 81;  Tag = "sv";  Act = "Statement-opinion" }
    @{ Row = 84;  Tag = "aa";  Act = "Agree/Accept" }
    @{ Row = 88;  Tag = "sv";  Act = "Statement-opinion" }
    @{ Row = 92;  Tag = "sv";  Act = "Statement-opinion" }
    @{ Row = 93;  Tag = "sd";  Act = "Statement-non-opinion" }
    @{ Row = 96;  Tag = "ba";  Act = "Appreciation" }
    @{ Row = 97;  Tag = "%";   Act = "Uninterpretable" }
    @{ Row = 146; Tag = "sv";  Act = "Statement-opinion" }
    @{ Row = 151; Tag = "b";   Act = "Acknowledge (Backchannel)" }
    @{ Row = 161; Tag = "sv";  Act = "Statement-opinion" }
    @{ Row = 181; Tag = "sv";  Act = "Statement-opinion" }
    @{ Row = 185; Tag = "sd";  Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
